$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New timestamp values for Z2:Z79 (column 26), in row order
$newTimestamps = @(
    "2025-11-13T06:52:28.600294",
    "2025-11-13T06:52:28.601294",
    "2025-11-13T06:52:28.601294",
    "2025-11-13T06:52:28.601294",
    "2025-11-13T06:52:28.601294",
    "2025-11-13T06:52:28.601294",
    "2025-11-13T06:52:28.601294",
    "2025-11-13T06:52:28.601294",
    "2025-11-13T06:52:28.601294",
    "2025-11-13T06:52:28.601294",
    "2025-11-13T06:52:28.601294",
    "2025-11-13T06:52:28.602297",
    "2025-11-13T06:52:28.602297",
    "2025-11-13T06:52:28.602297",
    "2025-11-13T06:52:28.602297",
    "2025-11-13T06:52:28.602297",
    "2025-11-13T06:52:28.602297",
    "2025-11-13T06:52:28.602297",
    "2025-11-13T06:52:28.602297",
    "2025-11-13T06:52:28.602297",
    "2025-11-13T06:52:28.602297",
    "2025-11-13T06:52:28.602297",
    "2025-11-13T06:52:28.603296",
    "2025-11-13T06:52:28.603296",
    "2025-11-13T06:52:28.603296",
    "2025-11-13T06:52:28.603296",
    "2025-11-13T06:52:28.603296",
    "2025-11-13T06:52:28.603296",
    "2025-11-13T06:52:28.603296",
    "2025-11-13T06:52:28.603296",
    "2025-11-13T06:52:28.603296",
    "2025-11-13T06:52:28.603296",
    "2025-11-13T06:52:28.603296",
    "2025-11-13T06:52:28.604295",
    "2025-11-13T06:52:28.604295",
    "2025-11-13T06:52:28.604295",
    "2025-11-13T06:52:28.604295",
    "2025-11-13T06:52:28.604295",
    "2025-11-13T06:52:28.604295",
    "2025-11-13T06:52:28.604295",
    "2025-11-13T06:52:28.604295",
    "2025-11-13T06:52:28.604295",
    "2025-11-13T06:52:28.604295",
    "2025-11-13T06:52:28.604295",
    "2025-11-13T06:52:28.605295",
    "2025-11-13T06:52:28.605295",
    "2025-11-13T06:52:28.605295",
    "2025-11-13T06:52:28.605295",
    "2025-11-13T06:52:28.605295",
    "2025-11-13T06:52:28.605295",
    "2025-11-13T06:52:28.605295",
    "2025-11-13T06:52:28.605295",
    "2025-11-13T06:52:28.605295",
    "2025-11-13T06:52:28.605295",
    "2025-11-13T06:52:28.608309",
    "2025-11-13T06:52:28.608309",
    "2025-11-13T06:52:28.893206",
    "2025-11-13T06:52:28.894203",
    "2025-11-13T06:52:28.894203",
    "2025-11-13T06:52:28.894203",
    "2025-11-13T06:52:28.894203",
    "2025-11-13T06:52:28.895204",
    "2025-11-13T06:52:28.895204",
    "2025-11-13T06:52:28.895204",
    "2025-11-13T06:52:28.895204",
    "2025-11-13T06:52:28.895204",
    "2025-11-13T06:52:28.896020",
    "2025-11-13T06:52:28.896020",
    "2025-11-13T06:52:28.896537",
    "2025-11-13T06:52:29.268200",
    "2025-11-13T06:52:29.268200",
    "2025-11-13T06:52:29.268200",
    "2025-11-13T06:52:29.268200",
    "2025-11-13T06:52:29.268200",
    "2025-11-13T06:52:29.268200",
    "2025-11-13T06:52:29.268200",
    "2025-11-13T06:52:29.268200",
    "2025-11-13T06:52:29.268200"
)

$startRow = 2
for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 26).Value = $newTimestamps[$i]
}
